$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each row per the scraped data refresh
$ws.Range("D2").Value = "26.397.86"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "1.833.76"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "254.18"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.5284"
$ws.Range("E7").Value = "  +0.56%  "
$ws.Range("D8").Value = "0.2837"
$ws.Range("E8").Value = "  -10.96%  "
$ws.Range("D9").Value = "0.06898"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "1.852.04"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").Value = "16.54"
$ws.Range("E11").Value = "  -11.26%  "
$ws.Range("D12").Value = "0.6962"
$ws.Range("E12").Value = "  -10.53%  "
$ws.Range("D13").Value = "0.07134"
$ws.Range("E13").Value = "  -7.64%  "
$ws.Range("D14").Value = "87.22"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "4.883"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "13.28"
$ws.Range("E18").Value = "  -3.62%  "
$ws.Range("D19").Value = "0.000007398"
$ws.Range("E19").Value = "  -6.42%  "
$ws.Range("D20").Value = "26.431.77"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("D21").Value = "2.086.46"
$ws.Range("E21").Value = "  +0.86%  "
$ws.Range("D22").Value = "4.515"
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("D23").Value = "5.843"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").Value = "8.995"
$ws.Range("E24").Value = "  -3.76%  "
$ws.Range("D25").Value = "141.98"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").Value = "1.673"
$ws.Range("E26").Value = "  +0.04%  "
$ws.Range("D27").Value = "2.041"
$ws.Range("E27").Value = "  -5.21%  "
$ws.Range("D28").Value = "16.58"
$ws.Range("D29").Value = "108.84"
$ws.Range("E29").Value = "  -2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.100"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("D31").Value = "0.08734"
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D32").Value = "3.873"
$ws.Range("E32").Value = "  -4.22%  "
$ws.Range("D33").Value = "0.04701"
$ws.Range("E33").Value = "  -2.68%  "
$ws.Range("D34").Value = "2.885"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.110"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("D36").Value = "0.7072"
$ws.Range("E36").Value = "  -2.24%  "
$ws.Range("D37").Value = "3.064"
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").Value = "2.187"
$ws.Range("E38").Value = "  -1.44%  "
$ws.Range("D39").Value = "0.01648"
$ws.Range("E39").Value = "  -6.34%  "
$ws.Range("D40").Value = "0.4484"
$ws.Range("E40").Value = "  -5.27%  "
$ws.Range("D41").Value = "0.8658"
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("D42").Value = "104.99"
$ws.Range("E42").Value = "  -3.90%  "
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").Value = "5.752"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("D45").Value = "7.061"
$ws.Range("E45").Value = "  -6.98%  "
$ws.Range("D46").Value = "8.704"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("D47").Value = "0.1192"
$ws.Range("E47").Value = "  -2.33%  "
$ws.Range("D50").Value = "0.05583"
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").Value = "0.8632"
$ws.Range("E51").Value = "  -3.13%  "

# Rows 48 and 49 swap coin identity (Aave <-> Elrond) along with refreshed values
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "33.54"
$ws.Range("E48").Value = "  -3.31%  "

$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "58.76"
$ws.Range("E49").Value = "  -0.98%  "
